# Update "想去人数" (column F) counts on the "展览" and "全部类型" sheets.
# Both sheets carry the same event listing, so the same row/value updates
# are applied to each.

$wb = $excel.ActiveWorkbook

$updates = @{
    5  = 138
    6  = 1342
    7  = 1576
    8  = 346
    9  = 443
    11 = 177
    14 = 120
    15 = 285
    16 = 322
    18 = 1771
    20 = 108
    22 = 697
    24 = 347
    25 = 4279
    27 = 292
    28 = 1125
    31 = 635
    33 = 326
    35 = 162
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
